# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.972.86'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '3.135.92'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.61'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.25'
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.125.49'
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  +6.70%  '
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.05'
$ws.Range("E14").Value = '  +4.50%  '
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("D16").Value = '3.641.25'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '63.735.08'
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").Value = '3.117.02'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.11'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.36'
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.30'
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.727'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.45'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.02'
$ws.Range("E24").Value = '  -2.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.37'
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.21'
$ws.Range("E26").Value = '  +2.50%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.17'
$ws.Range("E28").Value = '  +7.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.69'
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.03'
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.89'
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '0.0₃0874'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.42'
$ws.Range("E36").Value = '  +3.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.04'
$ws.Range("E37").Value = '  -0.86%  '
$ws.Range("E38").Value = '  -3.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.01'
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.40'
$ws.Range("E40").Value = '  -0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '443.56'
$ws.Range("E41").Value = '  +1.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.71'
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("D44").Value = '2.910.26'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("E46").Value = '  -2.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.22'
$ws.Range("E47").Value = '  +3.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.75'
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.44'
$ws.Range("E51").Value = '  -0.31%  '
